$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Highlight (red) the four TODO bullet items that gained a red highlight:
#    "Tutorial bit", "Bird Character", "What say", "Voice over"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Tutorial bit`r" -or $t -eq "Bird Character`r" -or $t -eq "What say`r" -or $t -eq "Voice over`r") {
        $p.Range.HighlightColorIndex = 6
    }
}

# ---------------------------------------------------------------------------
# 2) After the "Wincondition" bullet, add four new bullet items:
#    - Score in tiles pollunate en tijd maken          (en-US, no highlight)
#    - Extra uitleg text bij main screen toevoegen      (nl-NL, no highlight)
#    - Muziekje toevoegen                               (nl-NL, no highlight)
#    - KUNNEN VERLIEZEN                                 (nl-NL, green highlight)
# ---------------------------------------------------------------------------
$winCondition = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Wincondition`r") {
        $winCondition = $p
    }
}
$winIdx = $winCondition.Index

$winCondition.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($winIdx + 1)
$newPara1.Range.Text = "Score in tiles pollunate en tijd maken"
$newPara1.Range.HighlightColorIndex = 0
$newPara1.Range.LanguageID = "en-US"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($winIdx + 2)
$newPara2.Range.Text = "Extra uitleg text bij main screen toevoegen"
$newPara2.Range.HighlightColorIndex = 0
$newPara2.Range.LanguageID = "nl-NL"

$newPara2.Range.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item($winIdx + 3)
$newPara3.Range.Text = "Muziekje toevoegen"
$newPara3.Range.HighlightColorIndex = 0
$newPara3.Range.LanguageID = "nl-NL"

$newPara3.Range.InsertParagraphAfter()
$newPara4 = $d.Paragraphs.Item($winIdx + 4)
$newPara4.Range.Text = "KUNNEN VERLIEZEN"
$newPara4.Range.HighlightColorIndex = 4
$newPara4.Range.LanguageID = "nl-NL"

# ---------------------------------------------------------------------------
# 3) At the very end of the document, add two new (non-list) paragraphs:
#    - 16.00 start 17.00 Aanpassingen gemaakt, en je kan nu verliezen
#    - 17.00 start
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastIdx = $lastPara.Index

$lastPara.Range.InsertParagraphAfter()
$tailPara1 = $d.Paragraphs.Item($lastIdx + 1)
$tailPara1.Range.Text = "16.00 start 17.00 Aanpassingen gemaakt, en je kan nu verliezen"
$tailPara1.Range.HighlightColorIndex = 0
$tailPara1.Range.LanguageID = "nl-NL"

$tailPara1.Range.InsertParagraphAfter()
$tailPara2 = $d.Paragraphs.Item($lastIdx + 2)
$tailPara2.Range.Text = "17.00 start "
$tailPara2.Range.HighlightColorIndex = 0
$tailPara2.Range.LanguageID = "nl-NL"
